$d = $word.ActiveDocument
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -eq "") {
        # check its pPr-ish context via surrounding siblings; print only those with sz18/szCs18 style runs empty
    }
}
# Specifically search for paragraph index near where garbage collection ends (approx para index).
Write-Host "Total paragraphs: $($paras.Count)"
